# Insert a new row of data (row 4) into Sheet1, mirroring the "insert excel
# to sqlserver" commit: two new shared strings ("kaka", "testja") and a new
# data row A4:C4 = "kaka", 33, "testja".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "kaka"
$ws.Range("B4").Value = 33
$ws.Range("C4").Value = "testja"

# Update selection to match the post-edit state captured in the diff
# (<selection activeCell="C4" sqref="C4"/>).
$ws.Range("C4").Select()
